# foliar_inorganics.xlsx: swap the summary ANOVA table (Df/F/P + Residuals
# row) for a car::Anova-style Type-II table (Df/Chisq/Pr(>Chisq), no
# Residuals row) with refreshed model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so cells that no longer exist in the new table
# (the blank spacer cells in row 1/2 and the whole Residuals row) are not
# re-emitted with leftover formatting.
$ws.Cells.Clear()

# ---- Row 1: merged-looking header cells (every other column) ----
$ws.Range("C1").Value2 = "Foliar Ca"
$ws.Range("E1").Value2 = "Foliar P"
$ws.Range("G1").Value2 = "Foliar K"
$ws.Range("I1").Value2 = "Foliar Mg"
$ws.Range("K1").Value2 = "Foliar Al"
$ws.Range("M1").Value2 = "Foliar Zn"

# ---- Row 2: column sub-headers (Df once, then Chisq/Pr(>Chisq) pairs) ----
$ws.Range("B2").Value2 = "Df"
$ws.Range("C2").Value2 = "Chisq"
$ws.Range("D2").Value2 = "Pr(>Chisq)"
$ws.Range("E2").Value2 = "Chisq"
$ws.Range("F2").Value2 = "Pr(>Chisq)"
$ws.Range("G2").Value2 = "Chisq"
$ws.Range("H2").Value2 = "Pr(>Chisq)"
$ws.Range("I2").Value2 = "Chisq"
$ws.Range("J2").Value2 = "Pr(>Chisq)"
$ws.Range("K2").Value2 = "Chisq"
$ws.Range("L2").Value2 = "Pr(>Chisq)"
$ws.Range("M2").Value2 = "Chisq"
$ws.Range("N2").Value2 = "Pr(>Chisq)"

# ---- Row labels + Df column ----
$ws.Range("A3").Value2 = "Elevation"
$ws.Range("A4").Value2 = "Fire"
$ws.Range("A5").Value2 = "Elevation*Fire"
$ws.Range("B3").Value2 = 1
$ws.Range("B4").Value2 = 1
$ws.Range("B5").Value2 = 1

# ---- Data: refreshed Chisq / Pr(>Chisq) values per foliar nutrient ----
$ws.Range("C3").Value2 = 13.3023556745772
$ws.Range("D3").Value2 = 0.00026507287021465802
$ws.Range("E3").Value2 = 1.04839330274344
$ws.Range("F3").Value2 = 0.305877414032447
$ws.Range("G3").Value2 = 3.1584260394666899
$ws.Range("H3").Value2 = 0.075535978058448
$ws.Range("I3").Value2 = 0.017789007890601598
$ws.Range("J3").Value2 = 0.89389646820717195
$ws.Range("K3").Value2 = 0.34068541593773699
$ws.Range("L3").Value2 = 0.55943385381388
$ws.Range("M3").Value2 = 0.26682303049882899
$ws.Range("N3").Value2 = 0.60547091622308002

$ws.Range("C4").Value2 = 0.84346711131406205
$ws.Range("D4").Value2 = 0.35840705428702302
$ws.Range("E4").Value2 = 0.30897720231542197
$ws.Range("F4").Value2 = 0.57830849402261597
$ws.Range("G4").Value2 = 4.0713022193868502
$ws.Range("H4").Value2 = 0.043617636226494803
$ws.Range("I4").Value2 = 0.016031399198128599
$ws.Range("J4").Value2 = 0.89924499614360898
$ws.Range("K4").Value2 = 0.0206622884549072
$ws.Range("L4").Value2 = 0.88570275981497404
$ws.Range("M4").Value2 = 0.00409532076440046
$ws.Range("N4").Value2 = 0.94897445231980004

$ws.Range("C5").Value2 = 0.087781093423662998
$ws.Range("D5").Value2 = 0.76701721838815495
$ws.Range("E5").Value2 = 0.53514472949201197
$ws.Range("F5").Value2 = 0.46445190295720401
$ws.Range("G5").Value2 = 4.8631328939585696
$ws.Range("H5").Value2 = 0.027436471553291701
$ws.Range("I5").Value2 = 1.10892792619334
$ws.Range("J5").Value2 = 0.29231511234282198
$ws.Range("K5").Value2 = 0.187329294967573
$ws.Range("L5").Value2 = 0.66514878004127798
$ws.Range("M5").Value2 = 1.7938115815173299
$ws.Range("N5").Value2 = 0.180462448624684

# Df cells are whole numbers, the rest of the table uses the 0.000 format.
$ws.Range("B3:B5").NumberFormat = "General"
$ws.Range("C3:N5").NumberFormat = "0.000"

# Match the saved selection left behind by the edit (the data block).
$ws.Range("C3:N5").Select()

Write-Output "foliar_inorganics table refreshed"
